# Refresh the cryptocurrency price table: update the Price (D) and
# Volume(1h) (E) columns for every listed coin, and replace the
# BabyDogeCoin row (51) with the new EnergySwap entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.861.95"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "1.627.92"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "  -0.56%  "

$ws.Range("D5").Value = "'211.30"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.60%  "

$ws.Range("E8").Value = "  -0.75%  "

$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").Value = "1.858.57"
$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("D13").Value = "1.628.74"
$ws.Range("E13").Value = "  -1.11%  "

$ws.Range("E14").Value = "  -0.86%  "

$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").Value = "'65.17"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").Value = "27.849.24"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "'230.32"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("E20").Value = "  -1.90%  "

$ws.Range("D21").Value = "'0.994"
$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("D23").Value = "'10.31"
$ws.Range("E23").Value = "  -3.25%  "

$ws.Range("E24").Value = "  -4.23%  "

$ws.Range("D25").Value = "'154.19"

$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("E27").Value = "  -1.24%  "

$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("E30").Value = "  -1.49%  "

$ws.Range("E31").Value = "  -0.98%  "

$ws.Range("E32").Value = "  +1.95%  "

$ws.Range("D33").Value = "1.404.12"
$ws.Range("E33").Value = "  -2.46%  "

$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  +8.17%  "

$ws.Range("E37").Value = "  +0.77%  "

$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").Value = "'0.559"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").Value = "'0.866"
$ws.Range("E40").Value = "  -2.56%  "

$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("E42").Value = "  -0.55%  "

$ws.Range("E43").Value = "  -3.70%  "

$ws.Range("E44").Value = "  +1.40%  "

$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("D47").Value = "1.768.83"
$ws.Range("E47").Value = "  -1.15%  "

$ws.Range("D48").Value = "'87.74"
$ws.Range("E48").Value = "  -1.51%  "

$ws.Range("D49").Value = "'0.0997"
$ws.Range("E49").Value = "  -1.28%  "

$ws.Range("D50").Value = "'0.0506"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.46"
$ws.Range("E51").Value = "  -3.24%  "
